# Update the "Taxa de desocupação" table on Sheet1 for the new quarter:
#  - the quarter date moves from 01/07/2023 to 01/10/2023 for every region
#  - the ranking is refreshed, which reshuffles the region order (rows 2-7)
#  - the region that used to rank 5th ("Rio Grande do Norte") drops out of
#    the table entirely
#  - "Nordeste" (aggregate row) loses its ranking ("Colocação") value
#  - the trailing "Brasil" row shifts up, so the table now ends at row 9
#    instead of row 10 (dimension A1:E10 -> A1:E9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuarter = "01/10/2023"

# New data, in final row order (row 2 .. row 9).
# Columns: A=Região, D=Valor, E=Colocação ($null = no ranking cell for that row)
$rows = @(
    @{ A = "Amapá";         D = 14.2; E = "1º" },
    @{ A = "Bahia";         D = 12.7; E = "2º" },
    @{ A = "Pernambuco";    D = 11.9; E = "3º" },
    @{ A = "Sergipe";       D = 11.2; E = "4º" },
    @{ A = "Piauí";         D = 10.6; E = "5º" },
    @{ A = "Rio de Janeiro"; D = 10;  E = "6º" },
    @{ A = "Nordeste";      D = 10.4; E = $null },
    @{ A = "Brasil";        D = 7.4;  E = $null }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A

    # Keep the quarter column as literal text ("01/10/2023"), not an
    # auto-converted date serial: force text format, assign, then drop the
    # number-format override again so no stray style sticks to the cell.
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $newQuarter
    $cCell.ClearFormats()

    $ws.Cells.Item($r, 4).Value = $row.D

    if ($row.E) {
        $ws.Cells.Item($r, 5).Value = $row.E
    } else {
        $ws.Cells.Item($r, 5).ClearContents()
    }

    $r = $r + 1
}

# The previous table had 9 data rows (r2:r10); the refreshed one only has 8
# (r2:r9), so remove the now-unused trailing row entirely.
$ws.Rows.Item(10).Delete()
